$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to match repulled data / recalculated means
$ws.Range("F4").Value = -8
$ws.Range("F5").Value = -3
$ws.Range("F8").Value = -2
$ws.Range("F10").Value = -2
$ws.Range("F11").Value = -3
